$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = "'6"
$ws.Range("G2").Style = "Normal"

$ws.Range("D3").Value = "'21.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = "'6"
$ws.Range("G3").Style = "Normal"

$ws.Range("D4").Value = "'5.441"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = "'6"
$ws.Range("G4").Style = "Normal"

$ws.Range("D5").Value = "'0.05660"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'6"
$ws.Range("G5").Style = "Normal"

$ws.Range("D6").Value = "'3.379"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = "'6"
$ws.Range("G6").Style = "Normal"

$ws.Range("B7").Value = "'MXToken"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "'0.7998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6MXTokenMX"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'6"
$ws.Range("G7").Style = "Normal"

$ws.Range("B8").Value = "'FTXToken"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "'1.035"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'7FTXTokenFTT"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'6"
$ws.Range("G8").Style = "Normal"

$ws.Range("B9").Value = "'One"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.01155"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'8OneONEBestin24h"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'6"
$ws.Range("G9").Style = "Normal"

$ws.Range("D10").Value = "'0.1433"
$ws.Range("D10").Style = "Normal"
$ws.Range("G10").Value = "'6"
$ws.Range("G10").Style = "Normal"

$ws.Range("D11").Value = "'0.07246"
$ws.Range("D11").Style = "Normal"
$ws.Range("G11").Value = "'6"
$ws.Range("G11").Style = "Normal"

$ws.Range("D12").Value = "'0.03162"
$ws.Range("D12").Style = "Normal"
$ws.Range("G12").Value = "'6"
$ws.Range("G12").Style = "Normal"

$ws.Range("D13").Value = "'0.02950"
$ws.Range("D13").Style = "Normal"
$ws.Range("G13").Value = "'6"
$ws.Range("G13").Style = "Normal"

$ws.Range("G14").Value = "'6"
$ws.Range("G14").Style = "Normal"

$ws.Range("D15").Value = "'0.001658"
$ws.Range("D15").Style = "Normal"
$ws.Range("G15").Value = "'6"
$ws.Range("G15").Style = "Normal"

$ws.Range("D16").Value = "'3.214"
$ws.Range("D16").Style = "Normal"
$ws.Range("G16").Value = "'6"
$ws.Range("G16").Style = "Normal"

$ws.Range("D17").Value = "'0.04724"
$ws.Range("D17").Style = "Normal"
$ws.Range("G17").Value = "'6"
$ws.Range("G17").Style = "Normal"

$ws.Range("B18").Value = "'TigerCash"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.006395"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'17TigerCashTCH"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'6"
$ws.Range("G18").Style = "Normal"

$ws.Range("B19").Value = "'HotbitToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'0.005012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'18HotbitTokenHTB"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'6"
$ws.Range("G19").Style = "Normal"

$ws.Range("B20").Value = "'BitKan"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.001047"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'19BitKanKAN"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'6"
$ws.Range("G20").Style = "Normal"

$ws.Range("B21").Value = "'NitroEx"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'20NitroExNTX"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'6"
$ws.Range("G21").Style = "Normal"

$ws.Range("B22").Value = "'UpBots"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.0003203"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'21UpBotsUBXT"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'6"
$ws.Range("G22").Style = "Normal"

$ws.Range("B23").Value = "'LEO"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'3.873"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'22LEOLEO"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'6"
$ws.Range("G23").Style = "Normal"

$ws.Range("B24").Value = "'KuCoinToken"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'6.427"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'23KuCoinTokenKCS"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'6"
$ws.Range("G24").Style = "Normal"

$ws.Range("D25").Value = "'2.086"
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Value = "'6"
$ws.Range("G25").Style = "Normal"

$ws.Range("G26").Value = "'6"
$ws.Range("G26").Style = "Normal"

$ws.Range("D27").Value = "'0.1318"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = "'6"
$ws.Range("G27").Style = "Normal"

$ws.Range("G28").Value = "'6"
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").Value = "'6"
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").Value = "'6"
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").Value = "'6"
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").Value = "'6"
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").Value = "'6"
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").Value = "'6"
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").Value = "'6"
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").Value = "'6"
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").Value = "'6"
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").Value = "'6"
$ws.Range("G38").Style = "Normal"

$ws.Range("G39").Value = "'6"
$ws.Range("G39").Style = "Normal"

$ws.Range("D40").Value = "'0.04089"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'6"
$ws.Range("G40").Style = "Normal"

$ws.Range("B41").Value = "'KickToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.006926"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'6"
$ws.Range("G41").Style = "Normal"

$ws.Range("D42").Value = "'0.003503"
$ws.Range("D42").Style = "Normal"
$ws.Range("G42").Value = "'6"
$ws.Range("G42").Style = "Normal"

$ws.Range("B43").Value = "'BKEXToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'42BKEXTokenBKK"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'6"
$ws.Range("G43").Style = "Normal"

$ws.Range("D44").Value = "'0.008910"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = "'6"
$ws.Range("G44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005812"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = "'6"
$ws.Range("G45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("G46").Value = "'6"
$ws.Range("G46").Style = "Normal"

$ws.Range("D47").Value = "'0.7859"
$ws.Range("D47").Style = "Normal"
$ws.Range("G47").Value = "'6"
$ws.Range("G47").Style = "Normal"

$ws.Range("D48").Value = "'0.01537"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'6"
$ws.Range("G48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = "'6"
$ws.Range("G49").Style = "Normal"

$ws.Range("D50").Value = "'0.01011"
$ws.Range("D50").Style = "Normal"
$ws.Range("G50").Value = "'6"
$ws.Range("G50").Style = "Normal"

$ws.Range("G51").Value = "'6"
$ws.Range("G51").Style = "Normal"
